$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 187, shifting the existing rows 187-198 down to 188-199
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new weekly record
$ws.Cells.Item(187, 1).Value = 4
$ws.Cells.Item(187, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(187, 3).Value = "Los Lagos"
$ws.Cells.Item(187, 4).Value = 44610
$ws.Cells.Item(187, 5).Value = 10
$ws.Cells.Item(187, 6).Value = 100112032
$ws.Cells.Item(187, 7).Value = "Zapallo italiano"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 220
$ws.Cells.Item(187, 11).Value = 14000
$ws.Cells.Item(187, 12).Value = 15000
$ws.Cells.Item(187, 13).Value = 14455
$ws.Cells.Item(187, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(187, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(187, 16).Value = 289
$ws.Cells.Item(187, 17).Value = 50
$ws.Cells.Item(187, 18).Value = "Hortaliza"
